$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Fecha de ingreso" (entry date) values for rows 4-6 from
# 2025-12-16 (46007) to 2025-12-17 (46008)
$ws.Range("I4").Value = 46008
$ws.Range("I5").Value = 46008
$ws.Range("I6").Value = 46008

# Update the active selection to I5:I6 with I5 as the active cell
$ws.Range("I5:I6").Select()
